$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; existing columns A-D shift to B-E.
$ws.Columns("A").Insert()

# Populate the new first column ("TabName" / "CasesTab").
$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "CasesTab"

# New column C (the long query text, formerly column B) should wrap like
# column B already does.
$ws.Range("C2").WrapText = $true

# Resize columns to match the new layout.
$ws.Columns("A").ColumnWidth = 10
$ws.Columns("C").ColumnWidth = 127.67

# Update the sheet view: zoom out and move the selection.
$excel.ActiveWindow.Zoom = 40
$ws.Range("C11").Select()
